$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 134 (existing rows 134-182 shift down to 135-183)
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row 134 with the new record
$ws.Range("A134").Value = 4
$ws.Range("B134").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C134").Value = "Los Lagos"
$ws.Range("D134").Value = 44524
$ws.Range("E134").Value = 10
$ws.Range("F134").Value = 100112043
$ws.Range("G134").Value = "Pepino ensalada"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 70
$ws.Range("K134").Value = 11000
$ws.Range("L134").Value = 11000
$ws.Range("M134").Value = 11000
$ws.Range("N134").Value = "$/caja 60 unidades"
$ws.Range("O134").Value = "Región de Arica y Parinacota"
$ws.Range("P134").Value = 183
$ws.Range("Q134").Value = 60
$ws.Range("R134").Value = "Hortaliza"
